# rerun corona results with larger ds
# Updates the confidence/anchor-word statistics table (rows 3-27) with
# recomputed counts/percentages from the larger dataset, and appends a
# new word row (28, "please") that appears for the first time.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing changed cells (rows 3-27)
$ws.Range("B3").Value = 0.8529411764705882
$ws.Range("C3").Value = 29
$ws.Range("D3").Value = 29
$ws.Range("H3").Value = 5
$ws.Range("J3").Value = "best"
$ws.Range("K3").Value = 0.9322033898305084
$ws.Range("L3").Value = 55
$ws.Range("M3").Value = 55
$ws.Range("Q3").Value = 4
$ws.Range("B4").Value = 0.5787671232876712
$ws.Range("C4").Value = 169
$ws.Range("D4").Value = 169
$ws.Range("H4").Value = 123
$ws.Range("J4").Value = "love"
$ws.Range("K4").Value = 0.9130434782608695
$ws.Range("L4").Value = 42
$ws.Range("M4").Value = 42
$ws.Range("Q4").Value = 4
$ws.Range("B5").Value = 0.1782945736434109
$ws.Range("C5").Value = 92
$ws.Range("D5").Value = 92
$ws.Range("H5").Value = 424
$ws.Range("J5").Value = "interesting"
$ws.Range("K5").Value = 0.9090909090909091
$ws.Range("L5").Value = 30
$ws.Range("M5").Value = 30
$ws.Range("Q5").Value = 3
$ws.Range("B6").Value = 0.1587301587301587
$ws.Range("C6").Value = 30
$ws.Range("D6").Value = 30
$ws.Range("H6").Value = 159
$ws.Range("J6").Value = "great"
$ws.Range("K6").Value = 0.8303571428571429
$ws.Range("L6").Value = 93
$ws.Range("M6").Value = 93
$ws.Range("Q6").Value = 19
$ws.Range("J7").Value = "thanks"
$ws.Range("K7").Value = 0.8292682926829268
$ws.Range("L7").Value = 68
$ws.Range("M7").Value = 68
$ws.Range("Q7").Value = 14
$ws.Range("K8").Value = 0.7833333333333333
$ws.Range("L8").Value = 94
$ws.Range("M8").Value = 94
$ws.Range("Q8").Value = 26
$ws.Range("J9").Value = "special"
$ws.Range("K9").Value = 0.7777777777777778
$ws.Range("L9").Value = 28
$ws.Range("M9").Value = 28
$ws.Range("Q9").Value = 8
$ws.Range("J10").Value = "positive"
$ws.Range("K10").Value = 0.7586206896551724
$ws.Range("Q10").Value = 14
$ws.Range("J11").Value = "thank"
$ws.Range("K11").Value = 0.7421875
$ws.Range("L11").Value = 95
$ws.Range("M11").Value = 95
$ws.Range("Q11").Value = 33
$ws.Range("J12").Value = "safety"
$ws.Range("K12").Value = 0.7254901960784313
$ws.Range("L12").Value = 37
$ws.Range("M12").Value = 37
$ws.Range("Q12").Value = 14
$ws.Range("J13").Value = "confidence"
$ws.Range("K13").Value = 0.7222222222222222
$ws.Range("L13").Value = 26
$ws.Range("M13").Value = 26
$ws.Range("Q13").Value = 10
$ws.Range("J14").Value = "good"
$ws.Range("K14").Value = 0.71875
$ws.Range("L14").Value = 115
$ws.Range("M14").Value = 115
$ws.Range("Q14").Value = 45
$ws.Range("J15").Value = "safe"
$ws.Range("K15").Value = 0.6901408450704225
$ws.Range("L15").Value = 98
$ws.Range("M15").Value = 98
$ws.Range("Q15").Value = 44
$ws.Range("J16").Value = "support"
$ws.Range("K16").Value = 0.6886792452830188
$ws.Range("L16").Value = 73
$ws.Range("M16").Value = 73
$ws.Range("Q16").Value = 33
$ws.Range("J17").Value = "better"
$ws.Range("K17").Value = 0.6507936507936508
$ws.Range("L17").Value = 41
$ws.Range("M17").Value = 41
$ws.Range("Q17").Value = 22
$ws.Range("J18").Value = "well"
$ws.Range("K18").Value = 0.5957446808510638
$ws.Range("L18").Value = 56
$ws.Range("M18").Value = 56
$ws.Range("Q18").Value = 38
$ws.Range("J19").Value = "relief"
$ws.Range("K19").Value = 0.56
$ws.Range("L19").Value = 28
$ws.Range("M19").Value = 28
$ws.Range("Q19").Value = 22
$ws.Range("K20").Value = 0.5319148936170213
$ws.Range("L20").Value = 25
$ws.Range("M20").Value = 25
$ws.Range("Q20").Value = 22
$ws.Range("J21").Value = "fresh"
$ws.Range("K21").Value = 0.5208333333333334
$ws.Range("L21").Value = 25
$ws.Range("M21").Value = 25
$ws.Range("Q21").Value = 23
$ws.Range("J22").Value = "hand"
$ws.Range("K22").Value = 0.4830287206266319
$ws.Range("L22").Value = 185
$ws.Range("M22").Value = 185
$ws.Range("Q22").Value = 198
$ws.Range("J23").Value = "like"
$ws.Range("K23").Value = 0.4529411764705882
$ws.Range("L23").Value = 154
$ws.Range("M23").Value = 154
$ws.Range("Q23").Value = 186
$ws.Range("J24").Value = "care"
$ws.Range("K24").Value = 0.4382022471910113
$ws.Range("L24").Value = 39
$ws.Range("M24").Value = 39
$ws.Range("Q24").Value = 50
$ws.Range("J25").Value = "help"
$ws.Range("K25").Value = 0.423728813559322
$ws.Range("L25").Value = 125
$ws.Range("M25").Value = 125
$ws.Range("Q25").Value = 170
$ws.Range("J26").Value = "increase"
$ws.Range("K26").Value = 0.358974358974359
$ws.Range("L26").Value = 28
$ws.Range("M26").Value = 28
$ws.Range("Q26").Value = 50
$ws.Range("J27").Value = "protect"
$ws.Range("K27").Value = 0.3561643835616438
$ws.Range("L27").Value = 26
$ws.Range("M27").Value = 26
$ws.Range("Q27").Value = 47

# Add new row 28 (copy style/format from row 27 J:Q)
$ws.Range("J27:Q27").Copy()
$ws.Range("J28:Q28").PasteSpecial(-4122)

$ws.Range("J28").Value = "please"
$ws.Range("K28").Value = 0.3389121338912134
$ws.Range("L28").Value = 81
$ws.Range("M28").Value = 81
$ws.Range("N28").Value = 1
$ws.Range("O28").Value = 0
$ws.Range("P28").Value = $false
$ws.Range("Q28").Value = 158

$excel.CutCopyMode = $false